$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Persona" rows (rows 4-8), which shifts the "Denuncia" rows up.
$ws.Range("A4:H8").Delete()
